# Applies the "user registration form and routes!" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) New block describing the registration form living on the slideshow,
#    continuing the numbered list in column A (1 = Bootstrap, 2 = header,
#    3 = registration form on slideshow).
$ws.Range("A35").Value = 3
$ws.Range("B35").Value = "Registration form on slideshow"
$ws.Range("C35").Value = "registration component"

# 2) D11 previously held "user registration component"; it now describes the
#    registration component as a child of the slideshow component.
$ws.Range("D11").Value = "user registration component (child of slideshow component)"

# 3) New block (#4) documenting the app routes.
$ws.Range("A36").Value = 4
$ws.Range("B36").Value = "Implementing routes "
$ws.Range("C36").Value = "All possible routes without authenticatio and authorization!"
$ws.Range("D36").Value = "Home"
$ws.Range("D37").Value = "All Shops"
$ws.Range("D38").Value = "Add New Shop"
$ws.Range("D39").Value = "Shop Detail"
$ws.Range("D40").Value = "About Us"
$ws.Range("D41").Value = "Contact Us"

# Reflect the new viewport/selection the author ended up with after scrolling
# down to the newly-added rows.
$ws.Activate()
$ws.Range("A31").Select()
$excel.ActiveWindow.ScrollRow = 31
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D42:D43").Select()
